# Update sample bank and book Excel data:
# add a computed "amount" column (F) = debit (C) minus credit (D),
# with a header cell styled like the other headers (bold, centered,
# top-aligned) but boxed with a left/right border only, and leave the
# new column selected (matching the state Excel was left in when saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "amount" header in F1 -------------------------------------------
$hdr = $ws.Range("F1")
$hdr.Value = "amount"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.Item(7).LineStyle = 1   # xlEdgeLeft  -> xlContinuous
$hdr.Borders.Item(10).LineStyle = 1  # xlEdgeRight -> xlContinuous

# --- Formulas: F2 = C2-D2, F3:F10 share one formula (C-D) ----------------
$ws.Range("F2").Formula = "=C2-D2"
$ws.Range("F3:F10").Formula = "=C3-D3"

# --- Leave F2:F10 selected, active cell F2, matching the saved state -----
[void]$ws.Range("F2:F10").Select()
